$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old (current) row number -> new (target) row number for the
# data rows (2-95). This is a pure permutation: every "Beteckning" row
# keeps its full content (all columns A-Z) but moves to a new row index.
$rowMap = @{2=2; 3=3; 4=4; 5=5; 6=6; 7=7; 8=8; 9=9; 10=10; 11=11; 12=12; 13=13; 14=14; 15=17; 16=16; 17=15; 18=18; 19=20; 20=19; 21=23; 22=21; 23=22; 24=24; 25=25; 26=26; 27=27; 28=28; 29=29; 30=30; 31=31; 32=32; 33=38; 34=35; 35=41; 36=40; 37=33; 38=36; 39=39; 40=34; 41=42; 42=70; 43=37; 44=73; 45=89; 46=59; 47=92; 48=43; 49=91; 50=74; 51=67; 52=90; 53=54; 54=83; 55=88; 56=66; 57=56; 58=57; 59=49; 60=77; 61=58; 62=86; 63=61; 64=93; 65=55; 66=79; 67=62; 68=60; 69=51; 70=80; 71=87; 72=68; 73=53; 74=94; 75=75; 76=72; 77=78; 78=44; 79=69; 80=45; 81=64; 82=52; 83=63; 84=85; 85=71; 86=76; 87=95; 88=65; 89=48; 90=46; 91=47; 92=50; 93=82; 94=81; 95=84}

$lastCol = 26   # column Z
$firstRow = 2
$lastRow = 95

# Step 1: snapshot the full contents (A:Z) of every data row before we
# start overwriting anything, since the move is a full permutation and a
# later write could otherwise clobber data we still need to read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = New-Object 'object[]' $lastCol
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $rowVals[$c-1] = @{ F = $true; D = $cell.Formula }
        } else {
            $rowVals[$c-1] = @{ F = $false; D = $cell.Value2 }
        }
    }
    $snapshot[$r] = $rowVals
}

# Step 2: write every row back out at its new position, taken from the
# snapshot, so no data is lost regardless of write order.
foreach ($oldRow in $rowMap.Keys) {
    $newRow = $rowMap[$oldRow]
    $rowVals = $snapshot[$oldRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $info = $rowVals[$c-1]
        $cell = $ws.Cells.Item($newRow, $c)
        if ($info.D -eq $null) {
            $cell.Value2 = $null
        } elseif ($info.F) {
            $cell.Formula = $info.D
        } else {
            $cell.Value2 = $info.D
        }
    }
}

# Step 3: the "Förändrad" column (C) is refreshed to the new date serial
# for every data row, regardless of which original row its content came
# from.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 46066
}
